$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 542
$ws.Range("F5").Value = 77
$ws.Range("F6").Value = 2490
$ws.Range("F10").Value = 1512
$ws.Range("F12").Value = 610
$ws.Range("F13").Value = 1436
$ws.Range("F14").Value = 1436
$ws.Range("F15").Value = 1211
$ws.Range("F17").Value = 3542
$ws.Range("F18").Value = 640
$ws.Range("F19").Value = 3273
$ws.Range("F20").Value = 730
$ws.Range("F21").Value = 2115
$ws.Range("F23").Value = 283
$ws.Range("F24").Value = 9
$ws.Range("F25").Value = 1113
$ws.Range("F26").Value = 754
$ws.Range("F28").Value = 956
$ws.Range("F29").Value = 939
$ws.Range("F30").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 95
$ws.Range("F18").Value = 110
$ws.Range("F19").Value = 234
$ws.Range("F20").Value = 172
$ws.Range("F21").Value = 466

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 493
$ws.Range("G3").Value = "不可售"

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 542
$ws.Range("F10").Value = 77
$ws.Range("F11").Value = 493
$ws.Range("F12").Value = 2490
$ws.Range("F22").Value = 1512
$ws.Range("F25").Value = 1437
$ws.Range("F26").Value = 1437
$ws.Range("F27").Value = 95
$ws.Range("F29").Value = 1211
$ws.Range("F32").Value = 3542
$ws.Range("F33").Value = 640
$ws.Range("F34").Value = 3273
$ws.Range("F35").Value = 730
$ws.Range("F37").Value = 2115
$ws.Range("F39").Value = 283
$ws.Range("F40").Value = 1113
$ws.Range("F42").Value = 110
$ws.Range("F43").Value = 234
$ws.Range("F44").Value = 172
$ws.Range("F45").Value = 466
$ws.Range("F46").Value = 754
$ws.Range("F48").Value = 956
$ws.Range("F49").Value = 939
$ws.Range("F50").Value = 72

